$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in B12 from "Castilla-La Mancha" to "Castilla - La Mancha"
$ws.Range("B12").Value = "Castilla - La Mancha"

# The longer wrapped text now needs two lines, so the row grows taller
# (matches the diff's row 12 height change from 15 to 23.4)
$ws.Rows.Item(12).RowHeight = 23.4

# Update selection to B12 (matches the diff's sheetView selection change)
$ws.Range("B12").Select()
